$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 509.83334
$ws.Range("I11").Value = 509.83334
$ws.Range("K11").Value = 509.83334
$ws.Range("M11").Value = -369.83334
$ws.Range("H40").Value = 4080.7693
$ws.Range("I40").Value = 3150
$ws.Range("J40").Value = 4494.4443
$ws.Range("K40").Value = 3150
$ws.Range("L40").Value = 4494.4443
$ws.Range("M40").Value = -2975
$ws.Range("N40").Value = -4844.4443
$ws.Range("H70").Value = 2912.125
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 3113.8572
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 9341.571599999999
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -9881.571599999999
$ws.Range("H73").Value = 2912.125
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 3113.8572
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 9341.571599999999
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -11213.5716
$ws.Range("H74").Value = 10708.25
$ws.Range("I74").Value = 7699.8
$ws.Range("K74").Value = 7699.8
$ws.Range("M74").Value = -6763.8
$ws.Range("H77").Value = 10708.25
$ws.Range("I77").Value = 7699.8
$ws.Range("K77").Value = 38499
$ws.Range("M77").Value = -33819
$ws.Range("H98").Value = 90909736
$ws.Range("I98").Value = 111111680
$ws.Range("K98").Value = 111111680
$ws.Range("M98").Value = -111110182
$ws.Range("H122").Value = 90909736
$ws.Range("I122").Value = 111111680
$ws.Range("K122").Value = 333335040
$ws.Range("M122").Value = -333332590
$ws.Range("H132").Value = 4153.5625
$ws.Range("I132").Value = 3443.7144
$ws.Range("K132").Value = 10331.1432
$ws.Range("M132").Value = -7801.143199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1326.6666
$ws.Range("I2").Value = 992
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 992
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -879
$ws.Range("N2").Value = -3226
$ws.Range("H32").Value = 7054645.5
$ws.Range("I32").Value = 8071185.5
$ws.Range("K32").Value = 8071185.5
$ws.Range("M32").Value = -8070898.5
$ws.Range("H45").Value = 1903.6666
$ws.Range("I45").Value = 1297.6923
$ws.Range("K45").Value = 1297.6923
$ws.Range("M45").Value = -920.6922999999999
$ws.Range("H61").Value = 21793932
$ws.Range("I61").Value = 41672776
$ws.Range("J61").Value = 107921.45
$ws.Range("K61").Value = 41672776
$ws.Range("L61").Value = 107921.45
$ws.Range("M61").Value = -41672564
$ws.Range("N61").Value = -108345.45
$ws.Range("H116").Value = 1326.6666
$ws.Range("I116").Value = 992
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 992
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = 1302
$ws.Range("N116").Value = -7588
$ws.Range("H122").Value = 3387.7368
$ws.Range("I122").Value = 1821.2222
$ws.Range("K122").Value = 5463.6666
$ws.Range("M122").Value = -3013.6666
$ws.Range("H123").Value = 86197.60000000001
$ws.Range("J123").Value = 86197.60000000001
$ws.Range("L123").Value = 86197.60000000001
$ws.Range("N123").Value = -95997.60000000001
$ws.Range("H132").Value = 8628.191999999999
$ws.Range("I132").Value = 4307.4
$ws.Range("K132").Value = 12922.2
$ws.Range("M132").Value = -10392.2
$ws.Range("H136").Value = 21793932
$ws.Range("I136").Value = 41672776
$ws.Range("J136").Value = 107921.45
$ws.Range("K136").Value = 125018328
$ws.Range("L136").Value = 323764.35
$ws.Range("M136").Value = -125015778
$ws.Range("N136").Value = -328864.35

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1326.6666
$ws.Range("I3").Value = 992
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 992
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = -878
$ws.Range("N3").Value = -3228
$ws.Range("H22").Value = 1099.75
$ws.Range("I22").Value = 1099.75
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1099.75
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -926.75
$ws.Range("N22").ClearContents()
$ws.Range("H63").Value = 39333.332
$ws.Range("J63").Value = 39333.332
$ws.Range("L63").Value = 39333.332
$ws.Range("N63").Value = -40705.332
$ws.Range("H66").Value = 39333.332
$ws.Range("J66").Value = 39333.332
$ws.Range("L66").Value = 117999.996
$ws.Range("N66").Value = -124863.996
$ws.Range("H86").Value = 3253.3333
$ws.Range("J86").Value = 3130
$ws.Range("L86").Value = 3130
$ws.Range("N86").Value = -5376
$ws.Range("H89").Value = 3253.3333
$ws.Range("J89").Value = 3130
$ws.Range("L89").Value = 15650
$ws.Range("N89").Value = -26882
$ws.Range("H94").Value = 766.2632
$ws.Range("I94").Value = 769.6667
$ws.Range("K94").Value = 769.6667
$ws.Range("M94").Value = -318.6667
$ws.Range("H130").Value = 114248.625
$ws.Range("J130").Value = 114248.625
$ws.Range("L130").Value = 114248.625
$ws.Range("N130").Value = -124288.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 11109
$ws.Range("I16").Value = 11109
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 11109
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -10822
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 590807
$ws.Range("I31").Value = 11545.214
$ws.Range("J31").Value = 902717.1
$ws.Range("K31").Value = 11545.214
$ws.Range("L31").Value = 902717.1
$ws.Range("M31").Value = -11250.214
$ws.Range("N31").Value = -903307.1
$ws.Range("H34").Value = 590807
$ws.Range("I34").Value = 11545.214
$ws.Range("J34").Value = 902717.1
$ws.Range("K34").Value = 11545.214
$ws.Range("L34").Value = 902717.1
$ws.Range("M34").Value = -11343.214
$ws.Range("N34").Value = -903121.1
$ws.Range("H113").Value = 11109
$ws.Range("I113").Value = 11109
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 11109
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -8939
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2430.7646
$ws.Range("I132").Value = 2535.6667
$ws.Range("J132").Value = 2179
$ws.Range("K132").Value = 22821.0003
$ws.Range("L132").Value = 19611
$ws.Range("M132").Value = -20291.0003
$ws.Range("N132").Value = -24671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 261.13333
$ws.Range("I2").Value = 201.41667
$ws.Range("K2").Value = 201.41667
$ws.Range("M2").Value = -88.41667000000001
$ws.Range("H5").Value = 70005
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H102").Value = 3794.5715
$ws.Range("I102").Value = 3423.9167
$ws.Range("K102").Value = 3423.9167
$ws.Range("M102").Value = -1801.9167
$ws.Range("H126").Value = 4288.8
$ws.Range("I126").Value = 3966.5
$ws.Range("J126").Value = 4657.143
$ws.Range("K126").Value = 11899.5
$ws.Range("L126").Value = 13971.429
$ws.Range("M126").Value = -9429.5
$ws.Range("N126").Value = -18911.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2210.923
$ws.Range("I46").Value = 2406
$ws.Range("J46").Value = 1983.3334
$ws.Range("K46").Value = 2406
$ws.Range("L46").Value = 1983.3334
$ws.Range("M46").Value = -2218
$ws.Range("N46").Value = -2359.3334
$ws.Range("H68").Value = 2803.3572
$ws.Range("I68").Value = 2669.8
$ws.Range("K68").Value = 2669.8
$ws.Range("M68").Value = -1920.8
$ws.Range("H71").Value = 2803.3572
$ws.Range("I71").Value = 2669.8
$ws.Range("K71").Value = 13349
$ws.Range("M71").Value = -9605
$ws.Range("H93").Value = 111112990
$ws.Range("I93").Value = 142858130
$ws.Range("K93").Value = 142858130
$ws.Range("M93").Value = -142856882
$ws.Range("H100").Value = 8367.799999999999
$ws.Range("I100").Value = 9010.429
$ws.Range("J100").Value = 7805.5
$ws.Range("K100").Value = 9010.429
$ws.Range("L100").Value = 7805.5
$ws.Range("M100").Value = -8469.429
$ws.Range("N100").Value = -8887.5
$ws.Range("H136").Value = 95597.56
$ws.Range("I136").Value = 58485.555
$ws.Range("K136").Value = 175456.665
$ws.Range("M136").Value = -172906.665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 36618.89
$ws.Range("J15").Value = 80006.664
$ws.Range("L15").Value = 80006.664
$ws.Range("N15").Value = -80582.664
$ws.Range("H62").Value = 18188464
$ws.Range("I62").Value = 6273.75
$ws.Range("K62").Value = 6273.75
$ws.Range("M62").Value = -5649.75
$ws.Range("H65").Value = 18188464
$ws.Range("I65").Value = 6273.75
$ws.Range("K65").Value = 31368.75
$ws.Range("M65").Value = -28248.75
$ws.Range("H107").Value = 55556096
$ws.Range("J107").Value = 600.75
$ws.Range("L107").Value = 1802.25
$ws.Range("N107").Value = -5642.25
$ws.Range("H122").Value = 5657.85
$ws.Range("I122").Value = 3538.7058
$ws.Range("K122").Value = 10616.1174
$ws.Range("M122").Value = -8166.117400000001
$ws.Range("H132").Value = 272213.5
$ws.Range("I132").Value = 1997.0555
$ws.Range("K132").Value = 5991.166499999999
$ws.Range("M132").Value = -3461.166499999999
$ws.Range("H136").Value = 2239.6
$ws.Range("I136").Value = 601.7917
$ws.Range("K136").Value = 1805.3751
$ws.Range("M136").Value = 744.6249
